# "Fruta / hortaliza, semanal" — weekly refresh: insert the latest week's
# record for Femacal de La Calera / Pepino ensalada, pushing the existing
# history (rows 279-298) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 279, shifting rows 279:298
# down to 280:299 (same as Excel's Insert Sheet Rows on a selected row).
$ws.Rows.Item(279).Insert()

# Populate the newly inserted row 279 with this week's record. Columns
# A, B, C, E, F, G, H, I, N, O, Q, R repeat the values already used for
# this market/product/origin combination; D, J, K, L, M, P carry the new
# week's figures.
$ws.Range("A279").Value = 3
$ws.Range("B279").Value = "Femacal de La Calera"
$ws.Range("C279").Value = "Coquimbo"
$ws.Range("D279").Value = 44610
$ws.Range("E279").Value = 5
$ws.Range("F279").Value = 100112043
$ws.Range("G279").Value = "Pepino ensalada"
$ws.Range("H279").Value = "Sin especificar"
$ws.Range("I279").Value = "Primera"
$ws.Range("J279").Value = 90
$ws.Range("K279").Value = 15000
$ws.Range("L279").Value = 16000
$ws.Range("M279").Value = 15500
$ws.Range("N279").Value = '$/caja 70 unidades'
$ws.Range("O279").Value = "Limache"
$ws.Range("P279").Value = 221
$ws.Range("Q279").Value = 70
$ws.Range("R279").Value = "Hortaliza"
